$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.631.75"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "3.081.25"
$ws.Range("E3").Value = "  -2.49%  "
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").Value = "'590.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").Value = "'154.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.19%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").Value = "'0.541"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.82%  "
$ws.Range("D9").Value = "3.081.35"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("D10").Value = "'0.157"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.31%  "
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").Value = "'0.0000241"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("D14").Value = "'37.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").Value = "3.590.73"
$ws.Range("E15").Value = "  -2.51%  "
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "63.617.05"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'7.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("D19").Value = "3.078.11"
$ws.Range("E19").Value = "  -2.79%  "
$ws.Range("D20").Value = "'478.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.08%  "
$ws.Range("D21").Value = "'14.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").Value = "'0.711"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.37%  "
$ws.Range("D23").Value = "'7.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("D24").Value = "'2.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'81.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'12.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.92%  "
$ws.Range("D27").Value = "'10.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.95%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'7.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'2.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("E32").Value = "  -3.23%  "
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("D34").Value = "'27.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("D35").Value = "0.0₃0849"
$ws.Range("E35").Value = "  -1.83%  "
$ws.Range("D36").Value = "'3.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.11%  "
$ws.Range("E37").Value = "  -1.41%  "
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("D39").Value = "'2.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.99%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'50.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "'9.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("D42").Value = "'450.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.53%  "
$ws.Range("D43").Value = "'0.289"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("D44").Value = "'41.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.12%  "
$ws.Range("D45").Value = "'0.0362"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.55%  "
$ws.Range("E46").Value = "  +2.60%  "
$ws.Range("D47").Value = "2.813.31"
$ws.Range("E47").Value = "  -3.56%  "
$ws.Range("D48").Value = "'130.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("D49").Value = "'25.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.90%  "
$ws.Range("E51").Value = "  -0.05%  "
